$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Port-level landings: row 57 used to be the "Santa Cruz" / "MONTEREY AREA TOTALS"
# check row; it is now the sheet-wide "Totals" row instead.
$ws.Range("A57").Value = "MONTEREY AREA TOTALS"
$ws.Range("B57").Value = "Totals"

# Column A now holds the longer "MONTEREY AREA TOTALS" label, so re-fit its width.
[void]$ws.Columns("A").AutoFit()

# Select the whole of column A (mirrors clicking the column header).
[void]$ws.Columns("A").Select()
